$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h-volume-change (E) columns
# per the scraped data refresh. Values are written as text to
# preserve formatting (leading/trailing zeros, thousand-dot grouping,
# subscript-digit notation) exactly as scraped.

$ws.Range("D2").Value = '60.837.88'
$ws.Range("E2").Value = '  +3.51%  '

$ws.Range("D3").Value = '3.242.84'
$ws.Range("E3").Value = '  +2.40%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.17'
$ws.Range("E5").Value = '  +2.51%  '

$ws.Range("E6").Value = '  +4.55%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("E8").Value = '  -0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.38'
$ws.Range("E9").Value = '  +1.31%  '

$ws.Range("E10").Value = '  +2.62%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.433'
$ws.Range("E11").Value = '  -1.57%  '

$ws.Range("D12").Value = '3.809.69'
$ws.Range("E12").Value = '  +2.89%  '

$ws.Range("E13").Value = '  -2.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.30'
$ws.Range("E14").Value = '  +2.05%  '

$ws.Range("E15").Value = '  +2.57%  '

$ws.Range("D16").Value = '60.879.01'
$ws.Range("E16").Value = '  +3.74%  '

$ws.Range("D17").Value = '3.243.41'
$ws.Range("E17").Value = '  +2.93%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.30'
$ws.Range("E18").Value = '  +0.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.40'
$ws.Range("E19").Value = '  +3.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.39'
$ws.Range("E20").Value = '  +3.13%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.40'
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("E22").Value = '  -0.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.530'
$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.93'
$ws.Range("E24").Value = '  +0.39%  '

$ws.Range("E25").Value = '  +2.07%  '

$ws.Range("E26").Value = '  +3.14%  '

$ws.Range("D28").Value = '0.0₃0913'
$ws.Range("E28").Value = '  +6.23%  '

$ws.Range("E29").Value = '  +2.50%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.55'
$ws.Range("E30").Value = '  +0.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.18'
$ws.Range("E31").Value = '  +3.01%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.41'
$ws.Range("E32").Value = '  +4.85%  '

$ws.Range("E33").Value = '  +6.95%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.64'
$ws.Range("E34").Value = '  +4.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '158.96'
$ws.Range("E35").Value = '  +1.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.43'
$ws.Range("E36").Value = '  +6.97%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.39'
$ws.Range("E37").Value = '  +4.33%  '

$ws.Range("D38").Value = '2.798.45'
$ws.Range("E38").Value = '  +3.30%  '

$ws.Range("E39").Value = '  +3.89%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0314'
$ws.Range("E40").Value = '  +7.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.72'
$ws.Range("E41").Value = '  +1.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.28'
$ws.Range("E42").Value = '  -0.19%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.03'
$ws.Range("E43").Value = '  +2.37%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.730'
$ws.Range("E44").Value = '  +1.51%  '

$ws.Range("D45").Value = '3.289.27'
$ws.Range("E45").Value = '  +2.71%  '

$ws.Range("E46").Value = '  +2.42%  '

$ws.Range("E47").Value = '  +2.77%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.39'
$ws.Range("E48").Value = '  +6.42%  '

$ws.Range("E49").Value = '  +0.62%  '

$ws.Range("E50").Value = '  +7.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '275.53'
$ws.Range("E51").Value = '  +6.77%  '
